$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2085714285714286
$ws.Range("C2").Value = 0.5342857142857143
$ws.Range("J2").Value = 0.008571428571428572
$ws.Range("P2").Value = 0.1571428571428571
$ws.Range("S2").Value = 0.09142857142857143
$ws.Range("B3").Value = 0.01538461538461539
$ws.Range("C3").Value = 0.04615384615384616
$ws.Range("J3").Value = 0.05128205128205128
$ws.Range("P3").Value = 0.7487179487179487
$ws.Range("S3").Value = 0.1384615384615385
$ws.Range("J4").Value = 0.0303030303030303
$ws.Range("O4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("B6").Value = 0.07860262008733625
$ws.Range("D6").Value = 0.004366812227074236
$ws.Range("E6").Value = 0.004366812227074236
$ws.Range("F6").Value = 0.06986899563318777
$ws.Range("J6").Value = 0.3013100436681223
$ws.Range("O6").Value = 0.02620087336244541
$ws.Range("Q6").Value = 0.08733624454148471
$ws.Range("R6").Value = 0.07860262008733625
$ws.Range("S6").Value = 0.3493449781659388
$ws.Range("B7").Value = 0.1308411214953271
$ws.Range("D7").Value = 0.01869158878504673
$ws.Range("E7").Value = 0.004672897196261682
$ws.Range("F7").Value = 0.06542056074766354
$ws.Range("J7").Value = 0.1495327102803738
$ws.Range("O7").Value = 0.01869158878504673
$ws.Range("Q7").Value = 0.1261682242990654
$ws.Range("R7").Value = 0.06542056074766354
$ws.Range("S7").Value = 0.4205607476635514
$ws.Range("B8").Value = 0.1346578366445916
$ws.Range("D8").Value = 0.01324503311258278
$ws.Range("F8").Value = 0.08167770419426049
$ws.Range("J8").Value = 0.09933774834437085
$ws.Range("O8").Value = 0.01545253863134658
$ws.Range("Q8").Value = 0.1236203090507726
$ws.Range("R8").Value = 0.1125827814569536
$ws.Range("S8").Value = 0.4194260485651214
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.01333333333333333
$ws.Range("F9").Value = 0.05333333333333334
$ws.Range("J9").Value = 0.08
$ws.Range("O9").Value = 0.02666666666666667
$ws.Range("Q9").Value = 0.16
$ws.Range("R9").Value = 0.08666666666666667
$ws.Range("S9").Value = 0.48
$ws.Range("B10").Value = 0.1276962899050906
$ws.Range("D10").Value = 0.02243313201035375
$ws.Range("E10").Value = 0.001725625539257981
$ws.Range("F10").Value = 0.06125970664365833
$ws.Range("J10").Value = 0.1199309749784297
$ws.Range("O10").Value = 0.01207937877480587
$ws.Range("Q10").Value = 0.1570319240724763
$ws.Range("R10").Value = 0.07937877480586712
$ws.Range("S10").Value = 0.4184641932700604
$ws.Range("G11").Value = 0.130952380952381
$ws.Range("J11").Value = 0.09226190476190477
$ws.Range("K11").Value = 0.1904761904761905
$ws.Range("L11").Value = 0.5654761904761905
$ws.Range("S11").Value = 0.02083333333333333
$ws.Range("G12").Value = 0.7382198952879581
$ws.Range("J12").Value = 0.1989528795811518
$ws.Range("K12").Value = 0.005235602094240838
$ws.Range("L12").Value = 0.01047120418848168
$ws.Range("S12").Value = 0.04712041884816754
$ws.Range("G13").Value = 0.7333333333333333
$ws.Range("J13").Value = 0.2444444444444444
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("G14").Value = 0.4
$ws.Range("J14").Value = 0.2
$ws.Range("S14").Value = 0.4
$ws.Range("F15").Value = 0.025
$ws.Range("H15").Value = 0.15
$ws.Range("I15").Value = 0.05
$ws.Range("J15").Value = 0.445
$ws.Range("K15").Value = 0.05
$ws.Range("M15").Value = 0.015
$ws.Range("O15").Value = 0.055
$ws.Range("S15").Value = 0.21
$ws.Range("F16").Value = 0.04545454545454546
$ws.Range("H16").Value = 0.1136363636363636
$ws.Range("I16").Value = 0.05909090909090909
$ws.Range("J16").Value = 0.4045454545454545
$ws.Range("K16").Value = 0.1636363636363636
$ws.Range("M16").Value = 0.01818181818181818
$ws.Range("O16").Value = 0.06363636363636363
$ws.Range("S16").Value = 0.1318181818181818
$ws.Range("F17").Value = 0.0392156862745098
$ws.Range("H17").Value = 0.1568627450980392
$ws.Range("I17").Value = 0.08169934640522876
$ws.Range("J17").Value = 0.369281045751634
$ws.Range("K17").Value = 0.1437908496732026
$ws.Range("M17").Value = 0.009803921568627451
$ws.Range("O17").Value = 0.05228758169934641
$ws.Range("S17").Value = 0.1470588235294118
$ws.Range("F18").Value = 0.03141361256544502
$ws.Range("H18").Value = 0.2303664921465969
$ws.Range("I18").Value = 0.1099476439790576
$ws.Range("J18").Value = 0.3717277486910995
$ws.Range("K18").Value = 0.06282722513089005
$ws.Range("M18").Value = 0.03141361256544502
$ws.Range("N18").Value = 0.005235602094240838
$ws.Range("O18").Value = 0.05235602094240838
$ws.Range("S18").Value = 0.1047120418848168
$ws.Range("F19").Value = 0.01926040061633282
$ws.Range("H19").Value = 0.2326656394453005
$ws.Range("I19").Value = 0.06548536209553159
$ws.Range("J19").Value = 0.3305084745762712
$ws.Range("K19").Value = 0.1248073959938367
$ws.Range("M19").Value = 0.02311248073959938
$ws.Range("N19").Value = 0.002311248073959939
$ws.Range("O19").Value = 0.06779661016949153
$ws.Range("S19").Value = 0.1340523882896764
